# Update countries & provincias Spain
# Applies the COVID-19 stats refresh captured in the commit:
#  - updates the "Datos actualizados..." timestamp in A1
#  - refreshes case counts for several countries
#  - Australia's updated case count now overtakes Camerun, so those two rows swap
#  - Jamaica's updated case count now overtakes Santo Tome y Principe, so those two rows swap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 05:16"

function Set-CountryRow($Row, $Country, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Straightforward refreshes (ranking order unchanged)
Set-CountryRow 35  "Bolivia"       76789 1555 23582 50230 0 83 2977
Set-CountryRow 39  "Belgica"       68751 745  17546 41364 0 1  9841
Set-CountryRow 51  "Honduras"      42014 588  5554  35123 0 25 1337
Set-CountryRow 77  "Corea del Sur" 14336 31   13233 802   0 0  301
Set-CountryRow 91  "Haiti"         7424  12   4606  2657  0 0  161
Set-CountryRow 176 "Camboya"       239   5    164   75    0 0  0

# Australia overtakes Camerun -> rows 72/73 swap places
Set-CountryRow 72 "Australia" 17278 373 9983  7094 0 4 201
Set-CountryRow 73 "Camerun"   17255 0   15320 1544 0 0 391

# Jamaica overtakes Santo Tome y Principe -> rows 152/153 swap places
Set-CountryRow 152 "Jamaica"               878 14 726 142 0 0 10
Set-CountryRow 153 "Santo Tome y Principe" 871 0  778 78  0 0 15
